# "Generate Report for Handoff"
# Marks the ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md file as "Ready for handoff"
# (was "In Translation") across the Overview, zh-cn and de-de sheets, and
# stamps the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-21 12:15:42"

# zh-cn sheet: row 3 corresponds to ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-21 12:15:38"

# de-de sheet: row 3 corresponds to ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-21 12:15:42"
